# Update crypto price/volume snapshot cells (columns D and E, rows 2-51)
# with the latest scraped values from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.590.11"
$ws.Range("E2").Value = "  +3.97%  "
$ws.Range("D3").Value = "1.744.41"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'246.69"
$ws.Range("E5").Value = "  +3.85%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4815"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "'0.2695"
$ws.Range("E8").Value = "  +2.97%  "
$ws.Range("D9").Value = "'0.06263"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "1.744.53"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").Value = "'0.07128"
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").Value = "'15.84"
$ws.Range("E12").Value = "  +6.99%  "
$ws.Range("D13").Value = "'0.6227"
$ws.Range("E13").Value = "  +5.84%  "
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("D15").Value = "'77.57"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "'1.0000"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "26.581.17"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'0.000006904"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").Value = "'11.74"
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("D21").Value = "1.967.40"
$ws.Range("E21").Value = "  +4.28%  "
$ws.Range("D22").Value = "'4.641"
$ws.Range("E22").Value = "  +4.28%  "
$ws.Range("D23").Value = "'8.845"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "'5.375"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").Value = "'135.86"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "'15.40"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").Value = "'1.820"
$ws.Range("E27").Value = "  +5.64%  "
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("D29").Value = "'107.29"
$ws.Range("E29").Value = "  +2.51%  "
$ws.Range("D30").Value = "'4.015"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "'3.754"
$ws.Range("E31").Value = "  +3.48%  "
$ws.Range("D32").Value = "'0.07883"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "'0.04631"
$ws.Range("E33").Value = "  +7.97%  "
$ws.Range("D34").Value = "'2.618"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "'0.6463"
$ws.Range("E35").Value = "  +6.75%  "
$ws.Range("D36").Value = "'0.9989"
$ws.Range("E36").Value = "  +4.41%  "
$ws.Range("D37").Value = "'0.9453"
$ws.Range("E37").Value = "  +4.75%  "
$ws.Range("D38").Value = "'113.49"
$ws.Range("E38").Value = "  +17.53%  "
$ws.Range("D39").Value = "'1.996"
$ws.Range("E39").Value = "  +7.36%  "
$ws.Range("D40").Value = "'2.431"
$ws.Range("E40").Value = "  -6.15%  "
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "'5.773"
$ws.Range("E42").Value = "  +18.25%  "
$ws.Range("D43").Value = "'0.01514"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "'0.3926"
$ws.Range("E44").Value = "  +4.54%  "
$ws.Range("D45").Value = "'0.1219"
$ws.Range("E45").Value = "  +9.14%  "
$ws.Range("D46").Value = "'6.736"
$ws.Range("E46").Value = "  +8.37%  "
$ws.Range("D47").Value = "'0.05332"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").Value = "'7.981"
$ws.Range("E48").Value = "  +7.85%  "
$ws.Range("D49").Value = "'30.75"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").Value = "'1.274"
$ws.Range("D51").Value = "'0.3461"
$ws.Range("E51").Value = "  +3.77%  "
